# Auto-generated PowerShell COM-interop script
# Implements: 'added Malware detection slides'
$p = $ppt.ActivePresentation

# --- Slide 12: fill in the 'Malware Detection' overview content ---
$s12 = $p.Slides.Item(12)
$content12 = $s12.Shapes.Item(2)
$content12.TextFrame.TextRange.Text = "Two main methods for malware detection:`rPassive malware detection`rActive malware detection`rPassive method:`rA malicious attacker injects a malicious code into a user’s PC.`rActive method:`rThe malicious code collection system attempt to connect to a particular website and perform malicious action on the website in question. Called client honey pot."
$content12.TextFrame.TextRange.Paragraphs(2).IndentLevel = 2
$content12.TextFrame.TextRange.Paragraphs(3).IndentLevel = 2
$content12.TextFrame.TextRange.Paragraphs(5).IndentLevel = 2
$content12.TextFrame.TextRange.Paragraphs(7).IndentLevel = 2
$content12.TextFrame.AutoSize = 2

# --- New slide 13: client honey pot groups ---
$s13 = $p.Slides.Add(13, 2)
$title13 = $s13.Shapes.Item(1)
$title13.TextFrame.TextRange.Text = "Malware Detection"
$content13 = $s13.Shapes.Item(2)
$content13.TextFrame.TextRange.Text = "The client honey pot or active method divided into two groups:`rLow interaction client honey pot`rLow interaction client honey pot`rDetermined what is a malicious website`rThe actual website is not visited`rThe source code of target website is crawled`rComparing website source with the malicious action pattern of the system`rHigh interaction client honey pot`rVisit the website to check by using the web browser `rAnalyse the malicious website by monitoring malicious behaviour`rMonitor files, process creation, and registry modification"
$content13.TextFrame.TextRange.Paragraphs(1).ParagraphFormat.Bullet.Visible = $false
$content13.TextFrame.TextRange.Paragraphs(3).IndentLevel = 2
$content13.TextFrame.TextRange.Paragraphs(4).IndentLevel = 2
$content13.TextFrame.TextRange.Paragraphs(5).IndentLevel = 2
$content13.TextFrame.TextRange.Paragraphs(6).IndentLevel = 2
$content13.TextFrame.TextRange.Paragraphs(7).IndentLevel = 2
$content13.TextFrame.TextRange.Paragraphs(9).IndentLevel = 2
$content13.TextFrame.TextRange.Paragraphs(10).IndentLevel = 2
$content13.TextFrame.TextRange.Paragraphs(11).IndentLevel = 2
$content13.TextFrame.AutoSize = 2

# --- New slide 14: hybrid client honey pot ---
$s14 = $p.Slides.Add(14, 2)
$title14 = $s14.Shapes.Item(1)
$title14.TextFrame.TextRange.Text = "Malware Detection"
$content14 = $s14.Shapes.Item(2)
$content14.TextFrame.TextRange.Text = "There two possible solution for malware detection:`rHybrid client honey pot`rStudying structure of URLs words`rStudying the structure of URLs contains three steps:`rIdentify suspicious websites `rDerive lexical features for each suspicious websites`rPerform group analysis to pick out suspicious cluster "
$content14.TextFrame.TextRange.Paragraphs(1).ParagraphFormat.Bullet.Visible = $false
$content14.TextFrame.TextRange.Paragraphs(4).IndentLevel = 2
$content14.TextFrame.TextRange.Paragraphs(5).IndentLevel = 3
$content14.TextFrame.TextRange.Paragraphs(6).IndentLevel = 3
$content14.TextFrame.TextRange.Paragraphs(7).IndentLevel = 3
$content14.TextFrame.AutoSize = 2

Write-Host "Final slide count: $($p.Slides.Count)"

